$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.402677
$ws.Range("H2").Value = 37.20803100000001
$ws.Range("I2").Value = 0.4952943482020729
$ws.Range("J2").Value = 0.495294348202073
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 833.4308676666666
$ws.Range("N2").Value = 2500.292603
$ws.Range("O2").Value = 0.8518935545813505
$ws.Range("P2").Value = 0.8518935545813505
$ws.Range("Q2").Value = 10336.77385349941
$ws.Range("R2").Value = 93030.9646814947
$ws.Range("S2").Value = 0.4219380628539171
$ws.Range("T2").Value = 0.4219380628539171

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.402677
$ws.Range("H3").Value = 37.20803100000001
$ws.Range("I3").Value = 0.4952943482020729
$ws.Range("J3").Value = 0.495294348202073
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 17.73945766666667
$ws.Range("N3").Value = 53.218373
$ws.Range("O3").Value = 0.01813243333584592
$ws.Range("P3").Value = 0.01813243333584592
$ws.Range("Q3").Value = 220.0167635948404
$ws.Range("R3").Value = 1980.150872353563
$ws.Range("S3").Value = 0.008980891750395344
$ws.Range("T3").Value = 0.008980891750395344

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.402677
$ws.Range("H4").Value = 37.20803100000001
$ws.Range("I4").Value = 0.4952943482020729
$ws.Range("J4").Value = 0.495294348202073
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.344749666666667
$ws.Range("N4").Value = 4.034249
$ws.Range("O4").Value = 0.001374539410528448
$ws.Range("P4").Value = 0.001374539410528448
$ws.Range("Q4").Value = 16.67849576152434
$ws.Range("R4").Value = 150.106461853719
$ws.Range("S4").Value = 0.0006808016014157492
$ws.Range("T4").Value = 0.0006808016014157492

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.402677
$ws.Range("H5").Value = 37.20803100000001
$ws.Range("I5").Value = 0.4952943482020729
$ws.Range("J5").Value = 0.495294348202073
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 125.812397
$ws.Range("N5").Value = 377.437191
$ws.Range("O5").Value = 0.1285994726722751
$ws.Range("P5").Value = 0.1285994726722751
$ws.Range("Q5").Value = 1560.41052258677
$ws.Range("R5").Value = 14043.69470328092
$ws.Range("S5").Value = 0.06369459199634481
$ws.Range("T5").Value = 0.06369459199634481

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.161818666666666
$ws.Range("H6").Value = 15.485456
$ws.Range("I6").Value = 0.2061344991927113
$ws.Range("J6").Value = 0.2061344991927113
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 833.4308676666666
$ws.Range("N6").Value = 2500.292603
$ws.Range("O6").Value = 0.8518935545813505
$ws.Range("P6").Value = 0.8518935545813505
$ws.Range("Q6").Value = 4302.019010097995
$ws.Range("R6").Value = 38718.17109088197
$ws.Range("S6").Value = 0.1756046512391254
$ws.Range("T6").Value = 0.1756046512391254

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.161818666666666
$ws.Range("H7").Value = 15.485456
$ws.Range("I7").Value = 0.2061344991927113
$ws.Range("J7").Value = 0.2061344991927113
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 17.73945766666667
$ws.Range("N7").Value = 53.218373
$ws.Range("O7").Value = 0.01813243333584592
$ws.Range("P7").Value = 0.01813243333584592
$ws.Range("Q7").Value = 91.5678637203431
$ws.Range("R7").Value = 824.110773483088
$ws.Range("S7").Value = 0.003737720064829822
$ws.Range("T7").Value = 0.003737720064829823

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.161818666666666
$ws.Range("H8").Value = 15.485456
$ws.Range("I8").Value = 0.2061344991927113
$ws.Range("J8").Value = 0.2061344991927113
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.344749666666667
$ws.Range("N8").Value = 4.034249
$ws.Range("O8").Value = 0.001374539410528448
$ws.Range("P8").Value = 0.001374539410528448
$ws.Range("Q8").Value = 6.941353931393778
$ws.Range("R8").Value = 62.472185382544
$ws.Range("S8").Value = 0.0002833399930099263
$ws.Range("T8").Value = 0.0002833399930099263

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.161818666666666
$ws.Range("H9").Value = 15.485456
$ws.Range("I9").Value = 0.2061344991927113
$ws.Range("J9").Value = 0.2061344991927113
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 125.812397
$ws.Range("N9").Value = 377.437191
$ws.Range("O9").Value = 0.1285994726722751
$ws.Range("P9").Value = 0.1285994726722751
$ws.Range("Q9").Value = 649.4207793326774
$ws.Range("R9").Value = 5844.787013994096
$ws.Range("S9").Value = 0.0265087878957462
$ws.Range("T9").Value = 0.02650878789574621

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.073119
$ws.Range("H10").Value = 9.219357
$ws.Range("I10").Value = 0.1227233823836907
$ws.Range("J10").Value = 0.1227233823836907
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 833.4308676666666
$ws.Range("N10").Value = 2500.292603
$ws.Range("O10").Value = 0.8518935545813505
$ws.Range("P10").Value = 0.8518935545813505
$ws.Range("Q10").Value = 2561.232234612919
$ws.Range("R10").Value = 23051.09011151627
$ws.Range("S10").Value = 0.1045472584490886
$ws.Range("T10").Value = 0.1045472584490886

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.073119
$ws.Range("H11").Value = 9.219357
$ws.Range("I11").Value = 0.1227233823836907
$ws.Range("J11").Value = 0.1227233823836907
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.73945766666667
$ws.Range("N11").Value = 53.218373
$ws.Range("O11").Value = 0.01813243333584592
$ws.Range("P11").Value = 0.01813243333584592
$ws.Range("Q11").Value = 54.515464405129
$ws.Range("R11").Value = 490.639179646161
$ws.Range("S11").Value = 0.002225273549821799
$ws.Range("T11").Value = 0.002225273549821799

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.073119
$ws.Range("H12").Value = 9.219357
$ws.Range("I12").Value = 0.1227233823836907
$ws.Range("J12").Value = 0.1227233823836907
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.344749666666667
$ws.Range("N12").Value = 4.034249
$ws.Range("O12").Value = 0.001374539410528448
$ws.Range("P12").Value = 0.001374539410528448
$ws.Range("Q12").Value = 4.132575750877001
$ws.Range("R12").Value = 37.193181757893
$ws.Range("S12").Value = 0.0001686881256797356
$ws.Range("T12").Value = 0.0001686881256797356

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.073119
$ws.Range("H13").Value = 9.219357
$ws.Range("I13").Value = 0.1227233823836907
$ws.Range("J13").Value = 0.1227233823836907
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 125.812397
$ws.Range("N13").Value = 377.437191
$ws.Range("O13").Value = 0.1285994726722751
$ws.Range("P13").Value = 0.1285994726722751
$ws.Range("Q13").Value = 386.6364676562431
$ws.Range("R13").Value = 3479.728208906188
$ws.Range("S13").Value = 0.01578216225910061
$ws.Range("T13").Value = 0.01578216225910061

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.403408
$ws.Range("H14").Value = 13.210224
$ws.Range("I14").Value = 0.175847770221525
$ws.Range("J14").Value = 0.175847770221525
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 833.4308676666666
$ws.Range("N14").Value = 2500.292603
$ws.Range("O14").Value = 0.8518935545813505
$ws.Range("P14").Value = 0.8518935545813505
$ws.Range("Q14").Value = 3669.936150130341
$ws.Range("R14").Value = 33029.42535117307
$ws.Range("S14").Value = 0.1498035820392195
$ws.Range("T14").Value = 0.1498035820392195

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.403408
$ws.Range("H15").Value = 13.210224
$ws.Range("I15").Value = 0.175847770221525
$ws.Range("J15").Value = 0.175847770221525
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.73945766666667
$ws.Range("N15").Value = 53.218373
$ws.Range("O15").Value = 0.01813243333584592
$ws.Range("P15").Value = 0.01813243333584592
$ws.Range("Q15").Value = 78.11406980506133
$ws.Range("R15").Value = 703.026628245552
$ws.Range("S15").Value = 0.003188547970798954
$ws.Range("T15").Value = 0.003188547970798954

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.403408
$ws.Range("H16").Value = 13.210224
$ws.Range("I16").Value = 0.175847770221525
$ws.Range("J16").Value = 0.175847770221525
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.344749666666667
$ws.Range("N16").Value = 4.034249
$ws.Range("O16").Value = 0.001374539410528448
$ws.Range("P16").Value = 0.001374539410528448
$ws.Range("Q16").Value = 5.921481440197334
$ws.Range("R16").Value = 53.293332961776
$ws.Range("S16").Value = 0.000241709690423037
$ws.Range("T16").Value = 0.000241709690423037

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.403408
$ws.Range("H17").Value = 13.210224
$ws.Range("I17").Value = 0.175847770221525
$ws.Range("J17").Value = 0.175847770221525
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 125.812397
$ws.Range("N17").Value = 377.437191
$ws.Range("O17").Value = 0.1285994726722751
$ws.Range("P17").Value = 0.1285994726722751
$ws.Range("Q17").Value = 554.0033154489761
$ws.Range("R17").Value = 4986.029839040784
$ws.Range("S17").Value = 0.02261393052108352
$ws.Range("T17").Value = 0.02261393052108353
